# Apply Seraph_Profits Excel leve-profit value updates across all job sheets.
# Each sheet is a Table (Table_<JOB>) spanning A1:N141; columns H-N hold computed
# price/profit figures that were refreshed by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 250000100
$ws.Range("I5").Value = 333333440
$ws.Range("K5").Value = 333333440
$ws.Range("M5").Value = -333333325
$ws.Range("H9").Value = 293.14285
$ws.Range("I9").Value = 281.66666
$ws.Range("J9").Value = 362
$ws.Range("K9").Value = 281.66666
$ws.Range("L9").Value = 362
$ws.Range("M9").Value = -112.66666
$ws.Range("N9").Value = -700
$ws.Range("H15").Value = 1273.1628
$ws.Range("I15").Value = 1273.1628
$ws.Range("K15").Value = 3819.4884
$ws.Range("M15").Value = -3650.4884
$ws.Range("H55").Value = 1556.5
$ws.Range("I55").Value = 753.2222
$ws.Range("J55").Value = 3966.3333
$ws.Range("K55").Value = 753.2222
$ws.Range("L55").Value = 3966.3333
$ws.Range("M55").Value = -539.2222
$ws.Range("N55").Value = -4394.3333
$ws.Range("H101").Value = 796.8570999999999
$ws.Range("J101").Value = 999
$ws.Range("L101").Value = 2997
$ws.Range("N101").Value = -6241
$ws.Range("H137").Value = 1967.1111
$ws.Range("I137").Value = 2001
$ws.Range("K137").Value = 6003
$ws.Range("M137").Value = -3453

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8776.5
$ws.Range("I32").Value = 6004.8335
$ws.Range("K32").Value = 6004.8335
$ws.Range("M32").Value = -5717.8335
$ws.Range("H44").Value = 46328
$ws.Range("J44").Value = 46328
$ws.Range("L44").Value = 46328
$ws.Range("N44").Value = -47304
$ws.Range("H61").Value = 2516.6667
$ws.Range("I61").Value = 2516.6667
$ws.Range("K61").Value = 2516.6667
$ws.Range("M61").Value = -2304.6667
$ws.Range("H97").Value = 823.13336
$ws.Range("I97").Value = 897.0769
$ws.Range("K97").Value = 897.0769
$ws.Range("M97").Value = -401.0769
$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 1416
$ws.Range("J132").Value = 1941.6666
$ws.Range("L132").Value = 5824.9998
$ws.Range("N132").Value = -10884.9998
$ws.Range("H136").Value = 2516.6667
$ws.Range("I136").Value = 2516.6667
$ws.Range("K136").Value = 7550.000100000001
$ws.Range("M136").Value = -5000.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4821.5
$ws.Range("I105").Value = 4821.5
$ws.Range("K105").Value = 4821.5
$ws.Range("M105").Value = -3074.5
$ws.Range("H140").Value = 111112.336
$ws.Range("J140").Value = 111112.336
$ws.Range("L140").Value = 111112.336
$ws.Range("N140").Value = -121472.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 15000
$ws.Range("J26").Value = 15000
$ws.Range("L26").Value = 15000
$ws.Range("N26").Value = -15574
$ws.Range("H31").Value = 4418.6
$ws.Range("I31").Value = 3912.35
$ws.Range("K31").Value = 3912.35
$ws.Range("M31").Value = -3617.35
$ws.Range("H34").Value = 4418.6
$ws.Range("I34").Value = 3912.35
$ws.Range("K34").Value = 3912.35
$ws.Range("M34").Value = -3710.35
$ws.Range("H62").Value = 69459.164
$ws.Range("J62").Value = 102564
$ws.Range("L62").Value = 102564
$ws.Range("N62").Value = -103812
$ws.Range("H65").Value = 69459.164
$ws.Range("J65").Value = 102564
$ws.Range("L65").Value = 512820
$ws.Range("N65").Value = -519060
$ws.Range("H105").Value = 1100
$ws.Range("I105").Value = 1100
$ws.Range("K105").Value = 1100
$ws.Range("M105").Value = 647

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 202928
$ws.Range("J55").Value = 5000
$ws.Range("L55").Value = 15000
$ws.Range("N55").Value = -15354
$ws.Range("H128").Value = 2192667.8
$ws.Range("I128").Value = 2192667.8
$ws.Range("K128").Value = 6578003.399999999
$ws.Range("M128").Value = -6573023.399999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 656.25
$ws.Range("J2").Value = 1091.7142
$ws.Range("L2").Value = 1091.7142
$ws.Range("N2").Value = -1317.7142
$ws.Range("H3").Value = 66960.664
$ws.Range("I3").Value = 83500.836
$ws.Range("K3").Value = 83500.836
$ws.Range("M3").Value = -83384.836
$ws.Range("H97").Value = 1372.5555
$ws.Range("I97").Value = 1369.125
$ws.Range("K97").Value = 1369.125
$ws.Range("M97").Value = -873.125
$ws.Range("H109").Value = 27500
$ws.Range("J109").Value = 27500
$ws.Range("L109").Value = 27500
$ws.Range("N109").Value = -29580
$ws.Range("H122").Value = 61342.35
$ws.Range("I122").Value = 2204.6
$ws.Range("J122").Value = 145824.86
$ws.Range("K122").Value = 6613.799999999999
$ws.Range("L122").Value = 437474.58
$ws.Range("M122").Value = -4163.799999999999
$ws.Range("N122").Value = -442374.58
$ws.Range("H123").Value = 24700.215
$ws.Range("J123").Value = 24700.215
$ws.Range("L123").Value = 24700.215
$ws.Range("N123").Value = -29600.215
$ws.Range("H132").Value = 3251.6843
$ws.Range("I132").Value = 2839.1
$ws.Range("K132").Value = 8517.299999999999
$ws.Range("M132").Value = -5987.299999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2531.5625
$ws.Range("I40").Value = 2367
$ws.Range("K40").Value = 2367
$ws.Range("M40").Value = -2231
$ws.Range("H55").Value = 230.8
$ws.Range("I55").Value = 316.66666
$ws.Range("J55").Value = 102
$ws.Range("K55").Value = 316.66666
$ws.Range("L55").Value = 102
$ws.Range("M55").Value = -143.66666
$ws.Range("N55").Value = -448
$ws.Range("H93").Value = 1882.2
$ws.Range("I93").Value = 1882.2
$ws.Range("K93").Value = 1882.2
$ws.Range("M93").Value = -634.2
$ws.Range("H100").Value = 1882.3334
$ws.Range("I100").Value = 1882.3334
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1882.3334
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1341.3334
$ws.Range("N100").ClearContents()
$ws.Range("H132").Value = 6498.222
$ws.Range("I132").Value = 3789.8333
$ws.Range("J132").Value = 11915
$ws.Range("K132").Value = 11369.4999
$ws.Range("L132").Value = 35745
$ws.Range("M132").Value = -8839.499899999999
$ws.Range("N132").Value = -40805

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 933.4
$ws.Range("I100").Value = 1011.7222
$ws.Range("K100").Value = 2023.4444
$ws.Range("M100").Value = -1482.4444
$ws.Range("H126").Value = 2083.087
$ws.Range("I126").Value = 1408.5333
$ws.Range("J126").Value = 3347.875
$ws.Range("K126").Value = 4225.5999
$ws.Range("L126").Value = 10043.625
$ws.Range("M126").Value = -1755.5999
$ws.Range("N126").Value = -14983.625
